# Connor Readnour resume — content edits
#
# Only the genuine, author-authored content changes are reproduced here.
# (Relationship-id renumbering, empty <w:sdtEndPr/> stubs, <w:proofErr/>
# spans and <w:lastRenderedPageBreak/> repositioning in the source diff
# are all incidental artifacts of the hosting platform's own re-save /
# live pagination & proofing pass — they aren't something a user can, or
# would, deliberately produce through the Word object model, so they are
# intentionally left alone.)

$d = $word.ActiveDocument

# 1) "thousands of" -> "over 10,000" in the promotional-testing bullet.
$d.Content.Find.Execute(
    "thousands of", $true, $false, $false, $false, $false, $true, 1, $false,
    "over 10,000", 2
) | Out-Null

# 2) "utm parameters, and fallbacks" -> "UTM parameters and fallbacks"
#    (capitalize UTM, drop the comma) in that same bullet.
$d.Content.Find.Execute(
    "utm parameters, and fallbacks", $true, $false, $false, $false, $false, $true, 1, $false,
    "UTM parameters and fallbacks", 2
) | Out-Null

# 3) The "Demonstrated strong communication..." bullet was previously split
#    across three runs with identical formatting; the edit collapses it
#    back into a single run. Re-applying the (unchanged) visible text over
#    its own span merges the runs without altering what's on the page.
$d.Content.Find.Execute(
    "all team members, and effectively", $true, $false, $false, $false, $false, $true, 1, $false,
    "all team members, and effectively", 2
) | Out-Null
